$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Merge the "Alprojektvezető neve, aki felel ..." paragraph's split runs
#    (caused by spell-check proofErr markers) into a single run.
# ---------------------------------------------------------------------------
$rng = $d.Content.Duplicate
$rng.Find.Execute(
    "Alprojektvezető neve, aki felel az ösztöndíjas tevékenység végrehajtásáért: Dr. Gregorics Tibor",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Alprojektvezető neve, aki felel az ösztöndíjas tevékenység végrehajtásáért: Dr. Gregorics Tibor",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Merge the "Alprojektvezető aláírása" table-cell runs into a single run.
# ---------------------------------------------------------------------------
$rng = $d.Content.Duplicate
$rng.Find.Execute(
    "Alprojektvezető aláírása",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "Alprojektvezető aláírása",
    2) | Out-Null

# ---------------------------------------------------------------------------
# 3. Update the progress-report date from "2018. 10. 26" to "2018. 11. 26"
#    and relocate the "_GoBack" bookmark into the middle of the new text
#    (right after "2018. 11"), which also removes it from its old spot
#    after the "konferenciarészvétel..." text, since _GoBack is a
#    singleton bookmark.
# ---------------------------------------------------------------------------
$rng = $d.Content.Duplicate
$rng.Find.Execute(
    "2018. 10. 26",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "2018. 11. 26",
    2) | Out-Null

$rng2 = $d.Content.Duplicate
$rng2.Find.Execute("2018. 11. 26") | Out-Null
$bmPos = $rng2.Start + 8
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

Write-Output "done"
